$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) values for rows 2-6, columns A, D, E, G
# (B, C, F are empty inline strings in both before/after states).
# NOTE: use .Value2 (not .Value) -- .Value's getter does not resolve to the
# underlying scalar in this COM host and instead echoes the property's
# reflection signature.
$data = @{}
for ($r = 2; $r -le 6; $r++) {
    $data[$r] = @{
        A = $ws.Cells.Item($r, 1).Value2
        D = $ws.Cells.Item($r, 4).Value2
        E = $ws.Cells.Item($r, 5).Value2
        G = $ws.Cells.Item($r, 7).Value2
    }
}

# Mapping of source row (before) -> destination row (after), derived from the diff:
#   row 2 (Instrument Data) -> row 4
#   row 3 (Auxiliary Data)  -> row 6
#   row 4 (Validation)      -> row 5
#   row 5 (Verification)    -> row 3
#   row 6 (Uncertainty)     -> row 2
$map = @{
    2 = 4
    3 = 6
    4 = 5
    5 = 3
    6 = 2
}

foreach ($src in $map.Keys) {
    $dst = $map[$src]
    $ws.Cells.Item($dst, 1).Value2 = $data[$src].A
    $ws.Cells.Item($dst, 4).Value2 = $data[$src].D
    $ws.Cells.Item($dst, 5).Value2 = $data[$src].E
    $ws.Cells.Item($dst, 7).Value2 = $data[$src].G
}
